$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha2")

# Make row 24 ("Backlog - Principal") look like the other section-header rows
# (bold, larger font, taller row) by copying the style used by B2.
$ws.Range("B2").Copy()
$ws.Range("B24").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(24).RowHeight = 21

# Add the four new "working files" rows under the existing backlog list.
$ws.Range("B28").Value = "index.html (front)"
$ws.Range("C28").Value = "x"

$ws.Range("B29").Value = "style.css (front)"
$ws.Range("C29").Value = "x"

$ws.Range("B30").Value = "f1_main.js"
$ws.Range("C30").Value = "x"

$ws.Range("B31").Value = "f1_class.js"
$ws.Range("C31").Value = "x"

# Scroll / selection as captured in the saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("B34").Select()
